# Update the DESC column ("NO MAP") entries to the fuller note "NO MAP ON FILE"
# for the record-map / survey rows (rows 3-6 on the active sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = "NO MAP ON FILE"
$ws.Range("H4").Value = "NO MAP ON FILE"
$ws.Range("H5").Value = "NO MAP ON FILE"
$ws.Range("H6").Value = "NO MAP ON FILE"
